$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Ether
$ws.Range("H15").Value = 4215.2324
$ws.Range("I15").Value = 4215.2324
$ws.Range("K15").Value = 12645.6972
$ws.Range("M15").Value = -12476.6972

# Row 33: Clear Glass Lens
$ws.Range("H33").Value = 146.46666
$ws.Range("I33").Value = 116.69231
$ws.Range("K33").Value = 116.69231
$ws.Range("M33").Value = 112.30769

# Row 74: Wing Glue
$ws.Range("H74").Value = 2700
$ws.Range("I74").Value = 2700
$ws.Range("K74").Value = 2700
$ws.Range("M74").Value = -1764

# Row 77: Wing Glue
$ws.Range("H77").Value = 2700
$ws.Range("I77").Value = 2700
$ws.Range("K77").Value = 13500
$ws.Range("M77").Value = -8820

# Row 112: Superior Spiritbond Potion
$ws.Range("H112").Value = 2380.5144
$ws.Range("J112").Value = 2555.4194
$ws.Range("L112").Value = 7666.2582
$ws.Range("N112").Value = -9882.2582

# Row 129: Commanding Craftsman's Draught
$ws.Range("H129").Value = 888.5484
$ws.Range("I129").Value = 342.6
$ws.Range("J129").Value = 919.5682
$ws.Range("K129").Value = 1027.8
$ws.Range("L129").Value = 2758.7046
$ws.Range("M129").Value = 3972.2
$ws.Range("N129").Value = -12758.7046

# Row 132: Growth Formula Lambda
$ws.Range("H132").Value = 6064051
$ws.Range("I132").Value = 7579457.5
$ws.Range("K132").Value = 22738372.5
$ws.Range("M132").Value = -22735842.5

# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 1390.5714
$ws.Range("I137").Value = 940.4
$ws.Range("J137").Value = 1910
$ws.Range("K137").Value = 2821.2
$ws.Range("L137").Value = 5730
$ws.Range("M137").Value = -271.1999999999998
$ws.Range("N137").Value = -10830

# Row 141: Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 597.3889
$ws.Range("I141").Value = 597.3889
$ws.Range("K141").Value = 1792.1667
$ws.Range("M141").Value = 3387.8333

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Range("H32").Value = 4067.3447
$ws.Range("I32").Value = 3789.1875
$ws.Range("K32").Value = 3789.1875
$ws.Range("M32").Value = -3502.1875

# Row 132: Mountain Chromite Ingot
$ws.Range("H132").Value = 3546.3635
$ws.Range("I132").Value = 4653
$ws.Range("J132").Value = 2914
$ws.Range("K132").Value = 13959
$ws.Range("L132").Value = 8742
$ws.Range("M132").Value = -11429
$ws.Range("N132").Value = -13802

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Ash Lumber
$ws.Range("H16").Value = 45455716
$ws.Range("I16").Value = 62501036
$ws.Range("J16").Value = 1534.3334
$ws.Range("K16").Value = 62501036
$ws.Range("L16").Value = 1534.3334
$ws.Range("M16").Value = -62500749
$ws.Range("N16").Value = -2108.3334

# Row 31: Walnut Lumber
$ws.Range("H31").Value = 1706.32
$ws.Range("I31").Value = 1440.7693
$ws.Range("J31").Value = 1994
$ws.Range("K31").Value = 1440.7693
$ws.Range("L31").Value = 1994
$ws.Range("M31").Value = -1145.7693
$ws.Range("N31").Value = -2584

# Row 34: Walnut Lumber
$ws.Range("H34").Value = 1706.32
$ws.Range("I34").Value = 1440.7693
$ws.Range("J34").Value = 1994
$ws.Range("K34").Value = 1440.7693
$ws.Range("L34").Value = 1994
$ws.Range("M34").Value = -1238.7693
$ws.Range("N34").Value = -2398

# Row 99: Pine Lumber
$ws.Range("H99").Value = 1911.5
$ws.Range("I99").Value = 1898.8572
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1898.8572
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -400.8571999999999
$ws.Range("N99").Value = -4996

# Row 113: White Ash Lumber
$ws.Range("H113").Value = 45455716
$ws.Range("I113").Value = 62501036
$ws.Range("J113").Value = 1534.3334
$ws.Range("K113").Value = 62501036
$ws.Range("L113").Value = 1534.3334
$ws.Range("M113").Value = -62498866
$ws.Range("N113").Value = -5874.3334

# Row 122: Horse Chestnut Lumber
$ws.Range("H122").Value = 1055.6364
$ws.Range("I122").Value = 1014
$ws.Range("J122").Value = 1166.6666
$ws.Range("K122").Value = 3042
$ws.Range("L122").Value = 3499.9998
$ws.Range("M122").Value = -592
$ws.Range("N122").Value = -8399.9998

# Row 126: Red Pine Lumber
$ws.Range("H126").Value = 1911.5
$ws.Range("I126").Value = 1898.8572
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5696.571599999999
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3226.571599999999
$ws.Range("N126").Value = -10940

# Row 132: Ginseng Lumber
$ws.Range("H132").Value = 3888.9473
$ws.Range("I132").Value = 4466.815
$ws.Range("J132").Value = 2470.5454
$ws.Range("K132").Value = 13400.445
$ws.Range("L132").Value = 7411.6362
$ws.Range("M132").Value = -10870.445
$ws.Range("N132").Value = -12471.6362

# Row 134: Ceiba Lumber
$ws.Range("H134").Value = 3499.75
$ws.Range("I134").Value = 4719.6
$ws.Range("J134").Value = 1466.6666
$ws.Range("K134").Value = 14158.8
$ws.Range("L134").Value = 4399.9998
$ws.Range("M134").Value = -11623.8
$ws.Range("N134").Value = -9469.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 58: La Noscean Toast
$ws.Range("H58").Value = 2903.8333
$ws.Range("I58").Value = 2005
$ws.Range("J58").Value = 3083.6
$ws.Range("K58").Value = 6015
$ws.Range("L58").Value = 9250.799999999999
$ws.Range("M58").Value = -5887
$ws.Range("N58").Value = -9506.799999999999

# Row 131: Tsai tou Vounou
$ws.Range("H131").Value = 891.85565
$ws.Range("J131").Value = 932.34064
$ws.Range("L131").Value = 2797.02192
$ws.Range("N131").Value = -12877.02192

# Row 140: Mesquite Juice
$ws.Range("H140").Value = 38606.965
$ws.Range("I140").Value = 65322.688
$ws.Range("K140").Value = 195968.064
$ws.Range("M140").Value = -190788.064

$ws = $wb.Worksheets.Item("GSM")
# Row 97: Koppranickel Ingot
$ws.Range("H97").Value = 912.5
$ws.Range("I97").Value = 912.5
$ws.Range("K97").Value = 912.5
$ws.Range("M97").Value = -416.5

# Row 126: Phrygian Gold Ingot
$ws.Range("H126").Value = 1825.3334
$ws.Range("I126").Value = 1765.5555
$ws.Range("J126").Value = 2004.6666
$ws.Range("K126").Value = 5296.666499999999
$ws.Range("L126").Value = 6013.9998
$ws.Range("M126").Value = -2826.666499999999
$ws.Range("N126").Value = -10953.9998

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Leather
$ws.Range("H7").Value = 2168.4285
$ws.Range("I7").Value = 2162.875
$ws.Range("J7").Value = 2175.8333
$ws.Range("K7").Value = 2162.875
$ws.Range("L7").Value = 2175.8333
$ws.Range("M7").Value = -2050.875
$ws.Range("N7").Value = -2399.8333

# Row 16: Hard Leather
$ws.Range("H16").Value = 514.0476
$ws.Range("I16").Value = 534.75
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 534.75
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -364.75
$ws.Range("N16").Value = -440

# Row 61: Raptor Leather
$ws.Range("H61").Value = 2581.8
$ws.Range("I61").Value = 2102
$ws.Range("J61").Value = 2901.6667
$ws.Range("K61").Value = 2102
$ws.Range("L61").Value = 2901.6667
$ws.Range("M61").Value = -1900
$ws.Range("N61").Value = -3305.6667

# Row 108: Smilodonskin Trousers of Maiming
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null

# Row 113: Atrociraptor Leather
$ws.Range("H113").Value = 2581.8
$ws.Range("I113").Value = 2102
$ws.Range("J113").Value = 2901.6667
$ws.Range("K113").Value = 2102
$ws.Range("L113").Value = 2901.6667
$ws.Range("M113").Value = 68
$ws.Range("N113").Value = -7241.6667

# Row 126: Saiga Leather
$ws.Range("H126").Value = 2168.4285
$ws.Range("I126").Value = 2162.875
$ws.Range("J126").Value = 2175.8333
$ws.Range("K126").Value = 6488.625
$ws.Range("L126").Value = 6527.499899999999
$ws.Range("M126").Value = -4018.625
$ws.Range("N126").Value = -11467.4999

$ws = $wb.Worksheets.Item("WVR")
# Row 75: Ramie Turban of Crafting
$ws.Range("H75").Value = 35515
$ws.Range("J75").Value = 35515
$ws.Range("L75").Value = 35515
$ws.Range("N75").Value = -37387

# Row 78: Ramie Turban of Crafting
$ws.Range("H78").Value = 35515
$ws.Range("J78").Value = 35515
$ws.Range("L78").Value = 106545
$ws.Range("N78").Value = -115905

# Row 122: Dark Hempen Cloth
$ws.Range("H122").Value = 32501864
$ws.Range("I122").Value = 32501864
$ws.Range("K122").Value = 97505592
$ws.Range("M122").Value = -97503142
